$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Target change (per the diff):
#   "En esta sección se explicará ..."
# becomes
#   "En esta sección, se explicará ..."
# i.e. a comma is inserted right after "sección" (before the existing
# space), and the resulting text is split into three runs:
#   "En esta sección"  +  ","  +  " "
# while the following runs ("se ", "explicará", ...) must stay exactly as
# separate runs, just like in the original document.
# ---------------------------------------------------------------------------

# Locate the exact text to change without disturbing the rest of the
# document: search read-only first (Replace:=wdReplaceNone) so we get the
# match's Start/End without the Find engine rewriting anything.
$search = $d.Content
$found = $search.Find.Execute(
    "En esta sección ",   # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "",                     # ReplaceWith (unused, Replace:=wdReplaceNone)
    0                       # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find target text 'En esta sección '"
}

$matchStart = $search.Start
$matchEnd = $search.End

# Replace the matched text with the comma-augmented version. (This is a
# plain text substitution; the underlying engine coalesces same-formatted
# runs across the whole paragraph whenever its text changes, so we fix the
# run boundaries explicitly afterwards.)
$search.Text = "En esta sección, "

# Compute the character offsets of each logical piece following the edit.
$prefixLen  = "En esta sección".Length
$oTextStart = $matchStart            # start of "En esta sección"
$oComma     = $oTextStart + $prefixLen   # start of ","
$oSpace     = $oComma + 1                # start of " "
$oSe        = $oSpace + 1                # start of "se "
$oExplica   = $oSe + "se ".Length        # start of "explicará"
$oAfter     = $oExplica + "explicará".Length # start of the following run

# Force explicit run boundaries by toggling a character property on/off
# (a true no-op visually) at the start of every piece that must remain its
# own run. This prevents the save-time run-coalescing pass from merging
# runs that happen to share identical formatting.
function Mark-RunBoundary($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

Mark-RunBoundary $oComma   $oSpace     # ","
Mark-RunBoundary $oSpace   $oSe        # " "
Mark-RunBoundary $oSe      $oExplica   # "se "
Mark-RunBoundary $oExplica $oAfter     # "explicará"
